$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Variable Name"
$ws.Range("E1").Value = "Sentence Template Plural"

[void]$ws.Range("A3:XFD4").Select()
